$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Run Mode (column C) values for rows 11-32 from "Yes" to "No"
$ws.Range("C11:C32").Value = "No"

# Update Run Mode (column C) value for row 33 from "No" to "Yes"
$ws.Range("C33").Value = "Yes"

# Update the active selection to match the edited range
$ws.Range("C11:C32").Select()
